# Updates cryptos list prices/volume(1h) figures, and the two
# ranking swaps (Monero/LidoDAOToken at rows 32-33, FraxShare/Aave at rows 50-51),
# per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "46.701.03"
    "E2" = "  +2.85%  "
    "D3" = "2.263.65"
    "E3" = "  -0.71%  "
    "E4" = "  +0.05%  "
    "D5" = "299.94"
    "E5" = "  -0.33%  "
    "D6" = "99.85"
    "E6" = "  +4.55%  "
    "E7" = "  -1.37%  "
    "E8" = "  +0.14%  "
    "D9" = "0.507"
    "E9" = "  -1.38%  "
    "D10" = "35.10"
    "E10" = "  +2.37%  "
    "D11" = "0.0787"
    "E11" = "  -1.09%  "
    "D12" = "7.06"
    "E12" = "  -2.79%  "
    "E13" = "  -1.20%  "
    "D14" = "2.607.45"
    "E14" = "  -0.53%  "
    "D15" = "2.264.85"
    "E15" = "  -0.49%  "
    "D16" = "13.59"
    "E16" = "  -0.66%  "
    "D17" = "46.710.19"
    "E17" = "  +3.29%  "
    "D18" = "0.790"
    "E18" = "  -2.94%  "
    "D19" = "12.67"
    "E19" = "  -4.41%  "
    "D20" = "0.0₃0952"
    "E20" = "  +3.39%  "
    "D21" = "5.79"
    "E21" = "  -4.33%  "
    "D22" = "65.43"
    "E22" = "  -0.13%  "
    "D23" = "247.18"
    "E23" = "  +2.89%  "
    "D24" = "2.78"
    "E24" = "  -3.74%  "
    "E25" = "  -0.12%  "
    "E26" = "  -3.44%  "
    "D27" = "41.52"
    "E27" = "  +0.43%  "
    "D28" = "2.24"
    "E28" = "  -1.12%  "
    "D29" = "9.60"
    "E29" = "  -0.08%  "
    "D30" = "20.19"
    "E30" = "  +2.40%  "
    "D31" = "2.82"
    "E31" = "  +9.21%  "
    "B32" = "LidoDAOToken"
    "C32" = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
    "D32" = "3.32"
    "E32" = "  +12.85%  "
    "B33" = "Monero"
    "C33" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D33" = "145.85"
    "E33" = "  -3.77%  "
    "D34" = "5.35"
    "E34" = "  -3.66%  "
    "D35" = "0.0766"
    "E35" = "  -3.46%  "
    "D36" = "0.113"
    "E36" = "  +9.13%  "
    "E37" = "  -2.54%  "
    "D38" = "15.65"
    "E38" = "  +15.42%  "
    "D39" = "1.68"
    "E39" = "  -4.95%  "
    "D40" = "3.83"
    "E40" = "  -3.42%  "
    "D41" = "0.0295"
    "E41" = "  -5.61%  "
    "D42" = "3.09"
    "E42" = "  -4.78%  "
    "D43" = "0.999"
    "E43" = "  +0.03%  "
    "D44" = "1.784.17"
    "E44" = "  +0.88%  "
    "D45" = "90.76"
    "E45" = "  +18.61%  "
    "D46" = "1.89"
    "E46" = "  -3.94%  "
    "D47" = "70.97"
    "E47" = "  +0.88%  "
    "E48" = "  -4.76%  "
    "D49" = "4.79"
    "E49" = "  +1.00%  "
    "B50" = "Aave"
    "C50" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D50" = "93.98"
    "E50" = "  -1.83%  "
    "B51" = "FraxShare"
    "C51" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "D51" = "7.81"
    "E51" = "  -1.28%  "
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "0.999", "7.06")
    # are not silently coerced to the Number type by Excel -- the source
    # cells are plain text ("inlineStr") in the workbook.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    # Clear the format override back to the default style so no stray
    # cell formatting is introduced by the text coercion above.
    $cell.Style = "Normal"
}
